$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @(
    "https://www.youtube.com/watch?v=0rHUDWjR5gg&list=PL8dPuuaLjXtPAJr1ysd5yGIyiSFuh0mIL",
    "https://www.youtube.com/watch?v=L-Wtlev6suc&list=PL8dPuuaLjXtPAJr1ysd5yGIyiSFuh0mIL&index=2",
    "https://www.youtube.com/watch?v=01QWC-rZcfE&list=PL8dPuuaLjXtPAJr1ysd5yGIyiSFuh0mIL&index=3",
    "https://www.youtube.com/watch?v=AQ5vty8f9Xc&list=PL8dPuuaLjXtPAJr1ysd5yGIyiSFuh0mIL&index=4",
    "https://www.youtube.com/watch?v=PRgua7xceDA&index=5&list=PL8dPuuaLjXtPAJr1ysd5yGIyiSFuh0mIL",
    "https://www.youtube.com/watch?v=Yocja_N5s1I&list=PLBDA2E52FB1EF80C9",
    "https://www.youtube.com/watch?v=n7ndRwqJYDM&index=2&list=PLBDA2E52FB1EF80C9",
    "https://www.youtube.com/watch?v=sohXPx_XZ6Y&index=3&list=PLBDA2E52FB1EF80C9",
    "https://www.youtube.com/watch?v=Z3Wvw6BivVI&index=4&list=PLBDA2E52FB1EF80C9",
    "https://www.youtube.com/watch?v=Q-mkVSasZIM&list=PLBDA2E52FB1EF80C9&index=5",
    "https://www.youtube.com/watch?v=FSyAehMdpyI",
    "https://www.youtube.com/watch?v=hQpQ0hxVNTg",
    "https://www.youtube.com/watch?v=QiiyvzZBKT8",
    "https://www.youtube.com/watch?v=0RRVV4Diomg",
    "https://www.youtube.com/watch?v=rcKilE9CdaA",
    "https://www.youtube.com/watch?v=QnQe0xW_JY4",
    "https://www.youtube.com/watch?v=HVT3Y3_gHGg",
    "https://www.youtube.com/watch?v=H8WJ2KENlK0",
    "https://www.youtube.com/watch?v=cj8dDTHGJBY",
    "https://www.youtube.com/watch?v=dPKvHrD1eS4"
)

for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $urls[$i]
}

# Rows 21 and 22 mirror rows 1 and 2
$ws.Cells.Item(21, 2).Value = $urls[0]
$ws.Cells.Item(22, 2).Value = $urls[1]

$ws.Range("B22").Select()
